# Apply updated cryptocurrency price/volume figures (and the
# EnergySwap/Decentraland row swap) to match the refreshed data feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.260.31"
$ws.Range("E2").Value = "  -6.07%  "

$ws.Range("D3").Value = "1.842.44"
$ws.Range("E3").Value = "  -5.43%  "

$ws.Range("D4").Value = "`'1.004"
$ws.Range("E4").Value = "  -0.52%  "

$ws.Range("D5").Value = "`'332.84"
$ws.Range("E5").Value = "  +1.56%  "

$ws.Range("D6").Value = "`'1.005"
$ws.Range("E6").Value = "  -0.37%  "

$ws.Range("D7").Value = "`'0.4618"
$ws.Range("E7").Value = "  -4.89%  "

$ws.Range("D8").Value = "`'0.3882"
$ws.Range("E8").Value = "  -5.75%  "

$ws.Range("D9").Value = "`'46.09"
$ws.Range("E9").Value = "  -3.10%  "

$ws.Range("D10").Value = "`'0.07873"
$ws.Range("E10").Value = "  -4.28%  "

$ws.Range("D11").Value = "`'0.9694"
$ws.Range("E11").Value = "  -5.07%  "

$ws.Range("D12").Value = "`'22.01"
$ws.Range("E12").Value = "  -8.54%  "

$ws.Range("D13").Value = "1.867.12"
$ws.Range("E13").Value = "  -4.32%  "

$ws.Range("D14").Value = "`'5.787"
$ws.Range("E14").Value = "  -5.37%  "

$ws.Range("D15").Value = "`'6.944"
$ws.Range("E15").Value = "  -5.42%  "

$ws.Range("D16").Value = "`'0.06900"
$ws.Range("E16").Value = "  +0.65%  "

$ws.Range("D17").Value = "`'1.006"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").Value = "`'87.36"
$ws.Range("E18").Value = "  -5.01%  "

$ws.Range("D19").Value = "`'0.000009982"
$ws.Range("E19").Value = "  -4.11%  "

$ws.Range("D20").Value = "`'16.98"
$ws.Range("E20").Value = "  -5.58%  "

$ws.Range("D21").Value = "`'1.005"

$ws.Range("D22").Value = "28.292.92"
$ws.Range("E22").Value = "  -5.94%  "

$ws.Range("D23").Value = "`'5.368"
$ws.Range("E23").Value = "  -5.58%  "

$ws.Range("D24").Value = "`'11.12"
$ws.Range("E24").Value = "  -7.22%  "

$ws.Range("D25").Value = "`'2.169"
$ws.Range("E25").Value = "  -1.47%  "

$ws.Range("D26").Value = "2.110.43"
$ws.Range("E26").Value = "  -3.61%  "

$ws.Range("D27").Value = "`'153.72"
$ws.Range("E27").Value = "  -2.05%  "

$ws.Range("D28").Value = "`'19.34"
$ws.Range("E28").Value = "  -4.02%  "

$ws.Range("D29").Value = "`'5.909"
$ws.Range("E29").Value = "  -9.53%  "

$ws.Range("E30").Value = "  -6.42%  "

$ws.Range("D31").Value = "`'117.22"
$ws.Range("E31").Value = "  -3.35%  "

$ws.Range("D32").Value = "`'0.9492"
$ws.Range("E32").Value = "  -7.62%  "

$ws.Range("D33").Value = "`'0.09351"
$ws.Range("E33").Value = "  -3.09%  "

$ws.Range("D34").Value = "`'5.327"
$ws.Range("E34").Value = "  -5.69%  "

$ws.Range("D35").Value = "`'3.460"
$ws.Range("E35").Value = "  -2.91%  "

$ws.Range("E36").Value = "  -6.95%  "

$ws.Range("D37").Value = "`'0.06073"
$ws.Range("E37").Value = "  -6.91%  "

$ws.Range("E38").Value = "  -5.50%  "

$ws.Range("D39").Value = "`'1.154"
$ws.Range("E39").Value = "  -6.22%  "

$ws.Range("D40").Value = "`'1.004"
$ws.Range("E40").Value = "  -0.36%  "

$ws.Range("D41").Value = "`'7.621"
$ws.Range("E41").Value = "  -4.73%  "

$ws.Range("D42").Value = "`'0.5649"
$ws.Range("E42").Value = "  -5.60%  "

$ws.Range("D43").Value = "`'10.05"
$ws.Range("E43").Value = "  -6.93%  "

$ws.Range("E44").Value = "  -3.83%  "

$ws.Range("D45").Value = "`'2.386"

$ws.Range("D46").Value = "`'1.223"
$ws.Range("E46").Value = "  -4.57%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "`'11.74"
$ws.Range("E47").Value = "  -5.85%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "`'0.5332"
$ws.Range("E48").Value = "  -4.68%  "

$ws.Range("D49").Value = "`'0.07035"
$ws.Range("E49").Value = "  -6.89%  "

$ws.Range("D50").Value = "`'1.851"
$ws.Range("E50").Value = "  -7.40%  "

$ws.Range("D51").Value = "`'113.28"
$ws.Range("E51").Value = "  -4.11%  "
